$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 1 - "Subtitle 2" shape: update team name / author list
# -----------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(8)
$tr1 = $shp1.TextFrame.TextRange

# Step 1: "Corbin Getz, " (first run, chars 1-13) -> "Data Maniacs: Corbin Getz, "
$run1 = $tr1.Characters(1, 13)
$run1.Text = "Data Maniacs: Corbin Getz, "

# Step 2: split the remaining "Bilkis Khan, Robert stockwell" run into three
# runs ("Bilkis" / " Khan, Robert " / "stockwell") without altering their
# text or formatting, so the boundaries match the reviewed / spell-checked
# version of this line.
$full1 = $tr1.Text
$idxBilkis = $full1.IndexOf("Bilkis") + 1
$idxKhan = $full1.IndexOf(" Khan, Robert ") + 1
$idxStockwell = $full1.IndexOf("stockwell") + 1

$subBilkis = $tr1.Characters($idxBilkis, 6)
$subBilkis.Font.Bold = $subBilkis.Font.Bold

$subKhan = $tr1.Characters($idxKhan, 14)
$subKhan.Font.Bold = $subKhan.Font.Bold

$subStock = $tr1.Characters($idxStockwell, 9)
$subStock.Font.Bold = $subStock.Font.Bold

# -----------------------------------------------------------------
# Slide 14 - "Text Placeholder 3" shape: tidy punctuation
# -----------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$shp14 = $s14.Shapes.Item(3)
$tr14 = $shp14.TextFrame.TextRange
$tr14.Text = "The team was able to find some interesting things with this dataset. We answered some of the following questions. "
